# Atualização de bases das ligas, do dia: 19-04-2024 às 21:40
#
# This script reorders a handful of existing match rows (their underlying
# match records were re-sequenced upstream) and appends two brand new
# match rows at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: swap the full data payload (columns B:AC) of two rows, leaving
# column A (the sequential row id) untouched.
# ---------------------------------------------------------------------
function Swap-Rows($sheet, $rowA, $rowB) {
    $rangeA = "B" + $rowA + ":AC" + $rowA
    $rangeB = "B" + $rowB + ":AC" + $rowB
    $valsA = $sheet.Range($rangeA).Value2
    $valsB = $sheet.Range($rangeB).Value2
    $sheet.Range($rangeA).Value2 = $valsB
    $sheet.Range($rangeB).Value2 = $valsA
}

# ---------------------------------------------------------------------
# Helper: cyclically rotate the full data payload (columns B:AC) across
# three rows: row2 <- row1, row3 <- row2, row1 <- row3 (i.e. the new
# content of rowOrder[i] is the old content of rowOrder[i+1]).
# ---------------------------------------------------------------------
function Rotate-Rows($sheet, $row1, $row2, $row3) {
    $range1 = "B" + $row1 + ":AC" + $row1
    $range2 = "B" + $row2 + ":AC" + $row2
    $range3 = "B" + $row3 + ":AC" + $row3
    $vals1 = $sheet.Range($range1).Value2
    $vals2 = $sheet.Range($range2).Value2
    $vals3 = $sheet.Range($range3).Value2
    $sheet.Range($range1).Value2 = $vals2
    $sheet.Range($range2).Value2 = $vals3
    $sheet.Range($range3).Value2 = $vals1
}

# Bolivia Apertura block, 2024-01-23 fixtures: rows 27/28 swap
Swap-Rows $ws 27 28

# Bolivia Clausura block, 2024-05-14 fixtures: rows 107/108 swap
Swap-Rows $ws 107 108

# Bolivia Clausura block, 2024-05-31 fixtures: rows 143/144/145 rotate
Rotate-Rows $ws 143 144 145

# Bolivia Clausura block, 2024-06-03 fixtures: rows 148/149/150 rotate
Rotate-Rows $ws 148 149 150

# Bolivia Clausura block, 2024-06-04 fixtures: rows 153/154 swap
Swap-Rows $ws 153 154

# ---------------------------------------------------------------------
# Append two brand new match rows (224, 225) after the current last row
# (223). Clone the formatting of row 223 (bold-bordered id cell in A,
# custom date/time format in E) before filling in the new values so the
# generated styles line up with the rest of the sheet.
# ---------------------------------------------------------------------
$ws.Cells.Item(223, 1).Copy($ws.Cells.Item(224, 1))
$ws.Cells.Item(223, 5).Copy($ws.Cells.Item(224, 5))
$ws.Cells.Item(223, 1).Copy($ws.Cells.Item(225, 1))
$ws.Cells.Item(223, 5).Copy($ws.Cells.Item(225, 5))

# Row 224
$ws.Cells.Item(224, 1).Value2 = 222
$ws.Cells.Item(224, 2).Value2 = 8090748
$ws.Cells.Item(224, 3).Value2 = "Bolivia Primera División"
$ws.Cells.Item(224, 4).Value2 = "Bolivia Apertura"
$ws.Cells.Item(224, 5).Value2 = 45399.875
$ws.Cells.Item(224, 6).Value2 = "The Strongest"
$ws.Cells.Item(224, 7).Value2 = "San Jose de Oruro"
$ws.Cells.Item(224, 8).Value2 = 2
$ws.Cells.Item(224, 9).Value2 = 1
$ws.Cells.Item(224, 10).Value2 = "H"
$ws.Cells.Item(224, 11).Value2 = 1.533
$ws.Cells.Item(224, 12).Value2 = 4
$ws.Cells.Item(224, 13).Value2 = 5
$ws.Cells.Item(224, 14).Value2 = 1.285
$ws.Cells.Item(224, 15).Value2 = 5.5
$ws.Cells.Item(224, 16).Value2 = 11
$ws.Cells.Item(224, 17).Value2 = -1.75
$ws.Cells.Item(224, 18).Value2 = 1.975
$ws.Cells.Item(224, 19).Value2 = 1.825
$ws.Cells.Item(224, 20).Value2 = 3.25
$ws.Cells.Item(224, 21).Value2 = 2
$ws.Cells.Item(224, 22).Value2 = 1.8
$ws.Cells.Item(224, 23).Value2 = 0.2849999999999999
$ws.Cells.Item(224, 24).Value2 = -1
$ws.Cells.Item(224, 25).Value2 = -1
$ws.Cells.Item(224, 26).Value2 = -1
$ws.Cells.Item(224, 27).Value2 = 0.825
$ws.Cells.Item(224, 28).Value2 = -0.5
$ws.Cells.Item(224, 29).Value2 = 0.4

# Row 225
$ws.Cells.Item(225, 1).Value2 = 223
$ws.Cells.Item(225, 2).Value2 = 8090539
$ws.Cells.Item(225, 3).Value2 = "Bolivia Primera División"
$ws.Cells.Item(225, 4).Value2 = "Bolivia Apertura"
$ws.Cells.Item(225, 5).Value2 = 45400.83333333334
$ws.Cells.Item(225, 6).Value2 = "Bolivar"
$ws.Cells.Item(225, 7).Value2 = "San Antonio Bulo Bulo"
$ws.Cells.Item(225, 8).Value2 = 1
$ws.Cells.Item(225, 9).Value2 = 1
$ws.Cells.Item(225, 10).Value2 = "D"
$ws.Cells.Item(225, 11).Value2 = 1.363
$ws.Cells.Item(225, 12).Value2 = 4.5
$ws.Cells.Item(225, 13).Value2 = 7
$ws.Cells.Item(225, 14).Value2 = 1.125
$ws.Cells.Item(225, 15).Value2 = 8
$ws.Cells.Item(225, 16).Value2 = 19
$ws.Cells.Item(225, 17).Value2 = -2.25
$ws.Cells.Item(225, 18).Value2 = 1.8
$ws.Cells.Item(225, 19).Value2 = 2
$ws.Cells.Item(225, 20).Value2 = 3.5
$ws.Cells.Item(225, 21).Value2 = 1.75
$ws.Cells.Item(225, 22).Value2 = 2.05
$ws.Cells.Item(225, 23).Value2 = -1
$ws.Cells.Item(225, 24).Value2 = 7
$ws.Cells.Item(225, 25).Value2 = -1
$ws.Cells.Item(225, 26).Value2 = -1
$ws.Cells.Item(225, 27).Value2 = 1
$ws.Cells.Item(225, 28).Value2 = -1
$ws.Cells.Item(225, 29).Value2 = 1.05

Write-Output "Applied Bolivia Primera División update"
